# Fruta / hortaliza, semanal
# Insert a new week's pair of price rows (1a amarillo / 2a amarillo) right
# before row 816, shifting the existing rows 816:919 down to 818:921.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("816:817").Insert()

# Row 816 - new weekly entry, "1a amarillo"
$ws.Cells.Item(816, 1).Value = 4
$ws.Cells.Item(816, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(816, 3).Value = "Los Lagos"
$ws.Cells.Item(816, 4).Value = 45154
$ws.Cells.Item(816, 5).Value = 10
$ws.Cells.Item(816, 6).Value = "Fruta"
$ws.Cells.Item(816, 7).Value = 100102
$ws.Cells.Item(816, 8).Value = "Cítricos"
$ws.Cells.Item(816, 9).Value = 100102003
$ws.Cells.Item(816, 10).Value = "Limón"
$ws.Cells.Item(816, 11).Value = "Sin especificar"
$ws.Cells.Item(816, 12).Value = "1a amarillo"
$ws.Cells.Item(816, 13).Value = 600
$ws.Cells.Item(816, 14).Value = 11000
$ws.Cells.Item(816, 15).Value = 11000
$ws.Cells.Item(816, 16).Value = 11000
$ws.Cells.Item(816, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(816, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(816, 19).Value = 611
$ws.Cells.Item(816, 20).Value = 18

# Row 817 - new weekly entry, "2a amarillo"
$ws.Cells.Item(817, 1).Value = 4
$ws.Cells.Item(817, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(817, 3).Value = "Los Lagos"
$ws.Cells.Item(817, 4).Value = 45154
$ws.Cells.Item(817, 5).Value = 10
$ws.Cells.Item(817, 6).Value = "Fruta"
$ws.Cells.Item(817, 7).Value = 100102
$ws.Cells.Item(817, 8).Value = "Cítricos"
$ws.Cells.Item(817, 9).Value = 100102003
$ws.Cells.Item(817, 10).Value = "Limón"
$ws.Cells.Item(817, 11).Value = "Sin especificar"
$ws.Cells.Item(817, 12).Value = "2a amarillo"
$ws.Cells.Item(817, 13).Value = 600
$ws.Cells.Item(817, 14).Value = 9000
$ws.Cells.Item(817, 15).Value = 9000
$ws.Cells.Item(817, 16).Value = 9000
$ws.Cells.Item(817, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(817, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(817, 19).Value = 500
$ws.Cells.Item(817, 20).Value = 18

# Ensure the date cells keep the date number format used elsewhere in column D
$ws.Cells.Item(816, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(817, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
